# "Generate Report for Archive"
#
# The localization-status report is being regenerated: the zh-cn/de-de
# status moves from "Ready for handoff" to "In Translation", and the
# Status column(s) are narrowed to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Update the status text everywhere it appears --------------------------
# Overview!E2 (zh-cn status) and Overview!F2 (de-de status)
if ($overview.Range("E2").Value2 -eq $oldStatus) { $overview.Range("E2").Value2 = $newStatus }
if ($overview.Range("F2").Value2 -eq $oldStatus) { $overview.Range("F2").Value2 = $newStatus }

# zh-cn!C2 and de-de!C2 ("Status" column on the per-locale detail sheets)
if ($zhcn.Range("C2").Value2 -eq $oldStatus) { $zhcn.Range("C2").Value2 = $newStatus }
if ($dede.Range("C2").Value2 -eq $oldStatus) { $dede.Range("C2").Value2 = $newStatus }

# --- Narrow the Status columns to match the new, shorter text --------------
# Overview: columns E (zh-cn) and F (de-de)
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn / de-de: column C ("Status")
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
